# Mise a jour du fichier via Shiny
# Updates data values on sheets "pro", "ind", "conso" (sheet "VA" recalculates
# automatically since B = pro!B - conso!B), and restores the expected
# selection / active-cell state on each sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "pro": new yearly values (rows 2-26, column B)
# ------------------------------------------------------------------
$proValues = @(
    1071981.5531803041,
    1086535.7039443441,
    1104983.0594843852,
    963704.78306671581,
    945990.28082493623,
    930722.71517897362,
    1057337.6441696975,
    1117902.7684626032,
    1256432.5639295208,
    1189856.4708432767,
    1236989.1662027198,
    1302568.2104731307,
    1237643.039560745,
    1254727.2209089387,
    1206495.8700359438,
    1112020.8544587812,
    1462620.5966589863,
    1486238.4543876173,
    1605684.7941305023,
    1999783,
    2038236,
    2451701.3457828108,
    2600850.2537295958,
    3002201.8139454219,
    3415686.0456909742
)

$wsPro = $wb.Worksheets.Item("pro")
for ($i = 0; $i -lt $proValues.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $proValues[$i]
}

# ------------------------------------------------------------------
# Sheet "ind": new quarterly/monthly index values (rows 2-101, column B)
# ------------------------------------------------------------------
$indValues = @(
    265.98449413800546,
    276.88071335653399,
    258.20148041048526,
    310.15309704168322,
    269.59573018688764,
    280.63988592241498,
    261.70704751865401,
    314.36400432911421,
    274.17296429781499,
    285.40462925807674,
    266.15034646905679,
    319.70132047601857,
    291.75534650557063,
    301.43143953691617,
    246.09644318235729,
    308.47711733508066,
    285.65042593383009,
    275.24858969534341,
    243.34295738690787,
    273.44827188483606,
    284.49793268841631,
    282.2597364616152,
    255.71018935121799,
    345.25260858858434,
    316.40017755601588,
    303.72887909618174,
    271.47265178523509,
    346.76336461877412,
    332.92522008638144,
    271.1169144831598,
    263.00658557274818,
    307.15548368866047,
    347.69391476293589,
    316.2339719500356,
    263.5354723306304,
    317.89793408842741,
    298.99219535139093,
    306.02974438707474,
    312.92586924750651,
    326.27806659619881,
    330.72268267756579,
    353.68110743836326,
    305.26108164879599,
    398.08088043391416,
    339.81706407960678,
    361.679699290158,
    342.89711818412604,
    489.99504196870458,
    372.1771526161254,
    422.06192658086854,
    377.8967799278318,
    459.77776039015544,
    391.06267807680979,
    489.98572703593402,
    388.20451253232847,
    507.1766990358505,
    467.5499850413172,
    526.76306012261887,
    432.88169918273985,
    468.12690184229598,
    2.852781742627903,
    409.92066411204593,
    504.30393033864692,
    514.36723528285734,
    521.5125990724855,
    555.01791529404613,
    573.04375342028482,
    534.92208937365001,
    527.3200713469148,
    607.44176390556834,
    571.92872693838308,
    642.33055617763614,
    567.20429036640701,
    693.81396795998944,
    603.15222541343269,
    697.10662567405927,
    667.982190431395,
    781.1956097660468,
    692.96371495266874,
    776.52964767861329,
    802.93847719182349,
    864.16930171551667,
    737.81874296202341,
    768.92866218212748,
    1063.1755557769693,
    1069.3783541096388,
    995.63112308288385,
    942.28779751983234,
    847.95944819802912,
    1007.913630879233,
    833.43186754203577,
    939.0345863190181,
    905.54363725913663,
    1075.8630821461136,
    950.41021809405993,
    1008.2548201952303,
    984.2661460742205,
    877.71190050444682,
    1050.5122261470819,
    1139.0904303294642
)

$wsInd = $wb.Worksheets.Item("ind")
for ($i = 0; $i -lt $indValues.Length; $i++) {
    $wsInd.Cells.Item($i + 2, 2).Value = $indValues[$i]
}

# ------------------------------------------------------------------
# Sheet "conso": all values reset to 0 (rows 2-26, column B)
# ------------------------------------------------------------------
$wsConso = $wb.Worksheets.Item("conso")
for ($r = 2; $r -le 26; $r++) {
    $wsConso.Cells.Item($r, 2).Value = 0
}

# ------------------------------------------------------------------
# Selection / active cell updates on each sheet
# ("pro" is activated last so it remains the workbook's selected tab)
# ------------------------------------------------------------------
$wsVA = $wb.Worksheets.Item("VA")

$wsInd.Activate()
[void]$wsInd.Range("D16").Select()

$wsVA.Activate()
[void]$wsVA.Range("D16").Select()

$wsConso.Activate()
[void]$wsConso.Range("D16").Select()

$wsPro.Activate()
[void]$wsPro.Range("I2").Select()
